$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 5 new rows (one after each "method header" row), working from the
# bottom of the sheet upward so that not-yet-processed row numbers above
# each insertion point stay stable (same numbering as the original sheet).
$ws.Rows.Item(19).Insert()
$ws.Rows.Item(16).Insert()
$ws.Rows.Item(12).Insert()
$ws.Rows.Item(8).Insert()
$ws.Rows.Item(4).Insert()

# Fill in the newly inserted "properties" rows with their 3 values, and give
# them the same boxed-border + centered look as the other header/value rows
# (borders are applied cell-by-cell so each cell gets its own full box,
# matching the look of the existing rows such as B3:D3).
$newRows = @(4, 9, 14, 19, 23)
foreach ($r in $newRows) {
    $ws.Cells.Item($r, 2).Value = "properties"
    $ws.Cells.Item($r, 3).Value = "version"
    $ws.Cells.Item($r, 4).Value = "0.0.2"
    foreach ($c in @(2, 3, 4)) {
        $cell = $ws.Cells.Item($r, $c)
        $cell.Borders.LineStyle = 1
        $cell.HorizontalAlignment = -4108
    }
}

# Move the active selection to match the post-edit state.
$ws.Range("B23:D23").Select()
